$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a value into a cell and stamp it with the same
# number-format/style as a reference cell, without minting a brand new
# style record (Copy + PasteSpecial(xlPasteFormats) reuses the existing
# cellXfs entry, whereas setting .NumberFormat directly creates a new one).
# ---------------------------------------------------------------------
function Set-ValueWithStyleOf {
    param($targetCell, $value, $styleSourceCell)
    $targetCell.Value = $value
    $styleSourceCell.Copy() | Out-Null
    $targetCell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# Existing rows 6-15: "Waktu" (column D) is normalized to 16:00 (0.66..),
# and a new "Catatan" (column E) note is filled in for each bimbingan entry.
# ---------------------------------------------------------------------

# Row 5 - Bimbingan Seputar Judul (Waktu unchanged at 15:00)
$ws.Cells.Item(5, 5).Value = "Bimbingan Seputar Judul"

# Row 6 - Bimbingan Bab 1 (Waktu unchanged at 16:00)
$ws.Cells.Item(6, 5).Value = "Bimbingan Bab 1"

# Row 7 - Bimbingan Bab 1
$ws.Cells.Item(7, 4).Value = 0.66666666666666696
$ws.Cells.Item(7, 5).Value = "Bimbingan Bab 1"

# Row 8 - Bimbingan Rumusan Masalah
$ws.Cells.Item(8, 4).Value = 0.66666666666666696
$ws.Cells.Item(8, 5).Value = "Bimbingan Rumusan Masalah"

# Row 9 - Bimbingan Metodologi Masalah
$ws.Cells.Item(9, 4).Value = 0.66666666666666696
$ws.Cells.Item(9, 5).Value = "Bimbingan Metodologi Masalah"

# Row 10 - Bimbingan Tinjauan Pustaka
$ws.Cells.Item(10, 4).Value = 0.66666666666666696
$ws.Cells.Item(10, 5).Value = "Bimbingan Tinjauan Pustaka"

# Row 11 - Bimbingan Pengerjaan Project
$ws.Cells.Item(11, 4).Value = 0.66666666666666696
$ws.Cells.Item(11, 5).Value = "Bimbingan Pengerjaan Project"

# Row 12 - Bimbingan Pengerjaan Project
$ws.Cells.Item(12, 4).Value = 0.66666666666666696
$ws.Cells.Item(12, 5).Value = "Bimbingan Pengerjaan Project"

# Row 13 - Bimbingan Pengerjaan Project
$ws.Cells.Item(13, 4).Value = 0.66666666666666696
$ws.Cells.Item(13, 5).Value = "Bimbingan Pengerjaan Project"

# Row 14 - Bimbingan Pengerjaan Project
$ws.Cells.Item(14, 4).Value = 0.66666666666666696
$ws.Cells.Item(14, 5).Value = "Bimbingan Pengerjaan Project"

# Row 15 - Bimbingan Pengajuan HKI
$ws.Cells.Item(15, 4).Value = 0.66666666666666696
$ws.Cells.Item(15, 5).Value = "Bimbingan Pengajuan HKI"

# ---------------------------------------------------------------------
# New rows 16-18: Kelengkapan Sidang dan HKI entries
# ---------------------------------------------------------------------

# Row 16 - No 12
Set-ValueWithStyleOf $ws.Cells.Item(16, 3) 44566 $ws.Cells.Item(15, 3)
Set-ValueWithStyleOf $ws.Cells.Item(16, 4) 0.66666666666666696 $ws.Cells.Item(15, 4)
$ws.Cells.Item(16, 2).Value = 12
$ws.Cells.Item(16, 5).Value = "Bimbingan Seputar Pengganti Sidang"

# Row 17 - No 13
Set-ValueWithStyleOf $ws.Cells.Item(17, 3) 44567 $ws.Cells.Item(15, 3)
Set-ValueWithStyleOf $ws.Cells.Item(17, 4) 0.66666666666666696 $ws.Cells.Item(15, 4)
$ws.Cells.Item(17, 2).Value = 13
$ws.Cells.Item(17, 5).Value = "Bimbingan Seputar Pengganti Sidang"

# Row 18 - No 14
Set-ValueWithStyleOf $ws.Cells.Item(18, 3) 44576 $ws.Cells.Item(15, 3)
Set-ValueWithStyleOf $ws.Cells.Item(18, 4) 0.66666666666666696 $ws.Cells.Item(15, 4)
$ws.Cells.Item(18, 2).Value = 14
$ws.Cells.Item(18, 5).Value = "Bimbingan Seputar Pengganti Sidang"

# ---------------------------------------------------------------------
# Final selection, matching the post-edit cursor position in the workbook
# ---------------------------------------------------------------------
$ws.Range("I16").Select()
